# Append a new student row (13) to Sheet1: "sd", 10 — matching the
# existing bordered data-row style — then grow the dimension/AutoFilter
# range from A1:B12 to A1:B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row.
$ws.Range("A13").Value = "sd"
$ws.Range("B13").Value = 10

# Match the thin-border style used by the other data rows (A2:B12).
$ws.Range("A13:B13").Borders.LineStyle = 1

# Re-apply AutoFilter over the expanded range A1:B13 (drop the old
# A1:B12 filter first so the new call isn't treated as a toggle-off).
$ws.AutoFilterMode = $false
$ws.Range("A1:B13").AutoFilter()
